$d = $word.ActiveDocument

# 1. Remove the empty second paragraph entirely (its lone empty run is
#    dropped along with it -- paragraph 1 keeps its own run and now leads
#    straight into the table). Doing this before the bookmark surgery below
#    keeps bookmark range math simple.
$d.Paragraphs.Item(2).Range.Delete()

# 2. Rename the heading bookmark. Bookmark.Name is not actually persisted by
#    a direct assignment in this host, so re-create the bookmark at the same
#    (zero-width) spot under the new name and drop the old one. A zero-width
#    range sitting at the very start of the document confuses the host's
#    bookmark placement, so widen the helper range by one character (which
#    still yields two adjacent, content-free bookmark tags once it's added)
#    before handing it to Bookmarks.Add.
$bm = $d.Bookmarks.Item("_heading=h.q4l4bt5sniqq")
$bmRange = $d.Range($bm.Start, $bm.Start + 1)
$d.Bookmarks.Add("_xvh99mi2s8vh", $bmRange)
$bm.Delete()

# 3. First paragraph (Heading4): force explicit line spacing (276/auto) onto
#    this paragraph's own pPr (it already inherits 276/auto from docDefaults,
#    writing it explicitly is what the target XML does).
$p1 = $d.Paragraphs.Item(1)
$p1.Format.LineSpacingRule = 5   # wdLineSpaceMultiple
$p1.Format.LineSpacing = 13.8    # 276/240 * 12

# 4. Last paragraph (after the table): same explicit line-spacing fix.
$pLast = $d.Paragraphs.Last
$pLast.Format.LineSpacingRule = 5
$pLast.Format.LineSpacing = 13.8
